# Delete row 636 ("「幸運にも母親とその子２頭を５日間観察できた」") and shift
# all subsequent rows up by one, matching the target diff where the
# workbook's dimension shrinks from A1:C737 to A1:C736.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(636).Delete()
